# Auto-generated edit script applying the cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.443.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "'1.570.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").Value = "'288.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("D7").Value = "'0.3739"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.33%  "
$ws.Range("D8").Value = "'48.38"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.10%  "
$ws.Range("D9").Value = "'0.3331"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.07480"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "'1.132"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("D12").Value = "'1.004"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "'20.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.67%  "
$ws.Range("D14").Value = "'5.967"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("D15").Value = "'6.917"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").Value = "'1.571.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "'0.00001116"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").Value = "'88.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.45%  "
$ws.Range("D19").Value = "'0.06770"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "'6.377"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "'16.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'12.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("D24").Value = "'22.433.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "'2.393"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("D26").Value = "'2.563"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.40%  "
$ws.Range("D27").Value = "'153.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("D28").Value = "'19.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("D30").Value = "'124.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").Value = "'1.746.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").Value = "'1.055"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "'2.013"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'6.150"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.96%  "
$ws.Range("D35").Value = "'9.649"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.60%  "
$ws.Range("D36").Value = "'0.08279"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("D37").Value = "'0.02454"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("D38").Value = "'0.2268"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("D39").Value = "'0.06378"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.54%  "
$ws.Range("D40").Value = "'5.385"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("D41").Value = "'1.293"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.63%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.6293"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'11.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "'1.002"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'13.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("D46").Value = "'0.6152"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.86%  "
$ws.Range("D47").Value = "'3.781"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").Value = "'2.047"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("D49").Value = "'125.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("D50").Value = "'1.215"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.76%  "
$ws.Range("D51").Value = "'0.07264"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.50%  "
